$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts existing D:K data to F:M)
$ws.Range("D:E").EntireColumn.Insert(-4161)

# Copy the number formats/styles from the (now-shifted) old D:E columns (F:G)
# into the newly inserted blank D:E columns so they match the rest of the row.
$ws.Range("F:G").Copy()
$ws.Range("D:E").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 1715000
$ws.Cells.Item(8, 5).Value = 1712000
$ws.Cells.Item(9, 4).Value = 757000
$ws.Cells.Item(9, 5).Value = 502000
$ws.Cells.Item(10, 4).Value = 958000
$ws.Cells.Item(10, 5).Value = 1210000
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = "NA"
$ws.Cells.Item(15, 4).Value = 201000
$ws.Cells.Item(15, 5).Value = 201000
$ws.Cells.Item(17, 4).Value = 1049000
$ws.Cells.Item(17, 5).Value = 1040000
$ws.Cells.Item(18, 4).Value = 666000
$ws.Cells.Item(18, 5).Value = 672000
$ws.Cells.Item(20, 4).Value = -74000
$ws.Cells.Item(20, 5).Value = -19000
$ws.Cells.Item(21, 4).Value = 793000
$ws.Cells.Item(21, 5).Value = 854000
$ws.Cells.Item(22, 4).Value = 153000
$ws.Cells.Item(22, 5).Value = 134000
$ws.Cells.Item(23, 4).Value = 439000
$ws.Cells.Item(23, 5).Value = 519000
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(24, 5).Value = 3000
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 439000
$ws.Cells.Item(26, 5).Value = 516000
$ws.Cells.Item(27, 4).Value = 414000
$ws.Cells.Item(27, 5).Value = 491000
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = 74000
$ws.Cells.Item(32, 5).Value = 19000
$ws.Cells.Item(33, 4).Value = 414000
$ws.Cells.Item(33, 5).Value = 491000
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 414000
$ws.Cells.Item(35, 5).Value = 491000
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(41, 4).Value = 68000
$ws.Cells.Item(41, 5).Value = 37000
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(43, 4).Value = 706000
$ws.Cells.Item(43, 5).Value = 770000
$ws.Cells.Item(44, 4).Value = 77000
$ws.Cells.Item(44, 5).Value = 81000
$ws.Cells.Item(45, 4).Value = 46000
$ws.Cells.Item(45, 5).Value = 37000
$ws.Cells.Item(46, 4).Value = 897000
$ws.Cells.Item(46, 5).Value = 925000
$ws.Cells.Item(47, 4).Value = 4198000
$ws.Cells.Item(47, 5).Value = 4127000
$ws.Cells.Item(48, 4).Value = 14639000
$ws.Cells.Item(48, 5).Value = 14271000
$ws.Cells.Item(49, 4).Value = 3010000
$ws.Cells.Item(49, 5).Value = 3020000
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 35000
$ws.Cells.Item(52, 5).Value = 36000
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 22779000
$ws.Cells.Item(54, 5).Value = 22379000
$ws.Cells.Item(57, 4).Value = 365000
$ws.Cells.Item(57, 5).Value = 311000
$ws.Cells.Item(58, 4).Value = 0
$ws.Cells.Item(58, 5).Value = 0
$ws.Cells.Item(59, 4).Value = 821000
$ws.Cells.Item(59, 5).Value = 893000
$ws.Cells.Item(60, 4).Value = 1186000
$ws.Cells.Item(60, 5).Value = 1204000
$ws.Cells.Item(61, 4).Value = 14396000
$ws.Cells.Item(61, 5).Value = 13892000
$ws.Cells.Item(62, 4).Value = 333000
$ws.Cells.Item(62, 5).Value = 330000
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 16071000
$ws.Cells.Item(66, 5).Value = 15581000
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = 0
$ws.Cells.Item(72, 5).Value = 0
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 6708000
$ws.Cells.Item(76, 5).Value = 6798000
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = 414000
$ws.Cells.Item(81, 5).Value = 491000
$ws.Cells.Item(83, 4).Value = 201000
$ws.Cells.Item(83, 5).Value = 201000
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = 799000
$ws.Cells.Item(89, 5).Value = 737000
$ws.Cells.Item(91, 4).Value = -536000
$ws.Cells.Item(91, 5).Value = -521000
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = -659000
$ws.Cells.Item(94, 5).Value = -1073000
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(96, 5).Value = 0
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = -103000
$ws.Cells.Item(100, 5).Value = 366000
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(101, 5).Value = 0
$ws.Cells.Item(102, 4).Value = 37000
$ws.Cells.Item(102, 5).Value = 30000
